$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.401.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4458"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -6.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3839"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.41"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07806"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.015"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.41"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.849.97"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.833"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.084"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001024"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.05"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06499"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -9.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.463"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.387.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.76"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -7.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.260"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.057.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.30"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.036"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.451"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.88"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.475"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09303"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9226"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.590"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.212"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02211"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05927"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.282"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5898"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1847"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.247"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5656"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.359"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.918"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06841"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "107.73"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.57%  "
